# Applies the "Add image attribution, change font to Times New Roman" edit.
#
# Summary of changes (see xml_diff):
#   1. Slide 4 ("Second example: Chefs who need to interact"): the small
#      caption textbox under the opera-cake photo ("Google Shape;69;p15")
#      is resized taller and its text is replaced with a richer,
#      multi-run attribution string:
#         "Image: Opera Cake by Arnold Gatilao is licensed under CC BY 2.0"
#      where "Opera Cake", "Arnold Gatilao" and "CC BY 2.0" are hyperlinks,
#      and every run uses 12pt Times New Roman (latin/ea/cs).
#   2. Slide 8 ("Parallel output"): a leftover empty duplicate body
#      placeholder shape ("Google Shape;97;p19") is removed.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 4 - image attribution textbox
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)

# Locate the caption shape by name rather than a hard-coded index, in
# case shape ordering ever differs.
$caption = $null
for ($i = 1; $i -le $slide4.Shapes.Count; $i++) {
    $candidate = $slide4.Shapes.Item($i)
    if ($candidate.Name -eq "Google Shape;69;p15") {
        $caption = $candidate
    }
}
if ($caption -eq $null) {
    $caption = $slide4.Shapes.Item(4)
}

# Grow the textbox to fit the longer attribution text (EMU -> points).
$caption.Height = 792950 / 914400 * 72

$tr = $caption.TextFrame.TextRange
$tr.Text = "Image: Opera Cake by Arnold Gatilao is licensed under CC BY 2.0"

function Format-AttributionRun($range) {
    $range.Font.Size = 12
    $range.Font.Name = "Times New Roman"
    $range.Font.NameFarEast = "Times New Roman"
    $range.Font.NameComplexScript = "Times New Roman"
}

# "Image: "
$run1 = $tr.Characters(1, 7)
Format-AttributionRun($run1)

# "Opera Cake" (hyperlink to the photo)
$run2 = $tr.Characters(8, 10)
Format-AttributionRun($run2)
$run2.ActionSettings(1).Hyperlink.Address = "https://www.flickr.com/photos/arndog/1327813364"

# " by "
$run3 = $tr.Characters(18, 4)
Format-AttributionRun($run3)

# "Arnold Gatilao" (hyperlink to the author)
$run4 = $tr.Characters(22, 14)
Format-AttributionRun($run4)
$run4.ActionSettings(1).Hyperlink.Address = "https://www.flickr.com/photos/arndog/"

# " "
$run5 = $tr.Characters(36, 1)
Format-AttributionRun($run5)

# "is licensed under "
$run6 = $tr.Characters(37, 18)
Format-AttributionRun($run6)

# "CC BY 2.0" (hyperlink to the license)
$run7 = $tr.Characters(55, 9)
Format-AttributionRun($run7)
$run7.ActionSettings(1).Hyperlink.Address = "https://creativecommons.org/licenses/by/2.0/"

# ---------------------------------------------------------------------
# 2. Slide 8 - remove the leftover empty body placeholder shape
# ---------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)

for ($i = $slide8.Shapes.Count; $i -ge 1; $i--) {
    $shp = $slide8.Shapes.Item($i)
    if ($shp.Name -eq "Google Shape;97;p19") {
        $shp.Delete()
    }
}
